$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.182.63'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.59%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.431.33'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +1.41%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '549.14'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.62%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.91'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -3.21%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.641'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +6.20%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.627'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.07%  '
$ws.Range("E10").Value = '  +5.59%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.48'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -4.22%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000270'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.97%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.15'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.972.19'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.27%  '
$ws.Range("E15").Value = '  +0.59%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.421.55'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.22%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.25'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +1.43%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '65.232.93'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.70%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.79'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +2.98%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.980'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.78%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '412.79'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +5.62%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.00'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +4.80%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '85.18'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +2.44%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.11'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -2.92%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.72'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -11.63%  '
$ws.Range("E26").Value = '  +0.28%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.19'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +5.27%  '
$ws.Range("E28").Value = '  -1.69%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.87'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +4.42%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '29.71'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.39%  '
$ws.Range("E31").Value = '  -5.34%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '608.70'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -8.58%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.63'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.19%  '
$ws.Range("E34").Value = '  +0.13%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '58.98'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +2.18%  '
$ws.Range("E36").Value = '  -0.16%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.145'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +12.37%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '37.30'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.85%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0₃0773'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.85%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.377'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -5.70%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.176.00'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +4.78%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.29'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.53'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -10.16%  '
$ws.Range("B45").Value = 'ThetaToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.78'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.11%  '
$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.25'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +2.64%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0409'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.07%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.70'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.14%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.132'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +3.14%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '137.67'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.37%  '
$ws.Range("E51").Value = '  -2.40%  '
